$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "face/face102.png", "binden", "face")
    ,@(3, "car/car109.png", "parken", "car")
    ,@(4, "face/face095.png", "zielen", "face")
    ,@(5, "car/car086.png", "hassen", "car")
    ,@(6, "face/face066.png", "kriegen", "face")
    ,@(7, "car/car085.png", "hören", "car")
    ,@(8, "face/face076.png", "süßen", "face")
    ,@(9, "face/face068.png", "rechnen", "face")
    ,@(10, "car/car069.png", "fließen", "car")
    ,@(11, "car/car073.png", "proben", "car")
    ,@(12, "face/face080.png", "atmen", "face")
    ,@(13, "face/face100.png", "münzen", "face")
    ,@(14, "face/face089.png", "lügen", "face")
    ,@(15, "car/car087.png", "wecken", "car")
    ,@(16, "face/face074.png", "quellen", "face")
    ,@(17, "face/face104.png", "stoppen", "face")
    ,@(18, "face/face082.png", "duschen", "face")
    ,@(19, "car/car076.png", "lernen", "car")
    ,@(20, "car/car075.png", "kranken", "car")
    ,@(21, "car/car107.png", "danken", "car")
    ,@(22, "car/car097.png", "herrschen", "car")
    ,@(23, "face/face101.png", "streifen", "face")
    ,@(24, "face/face105.png", "dienen", "face")
    ,@(25, "face/face071.png", "spenden", "face")
    ,@(26, "car/car108.png", "nullen", "car")
    ,@(27, "car/car105.png", "achten", "car")
    ,@(28, "face/face092.png", "spüren", "face")
    ,@(29, "car/car083.png", "grenzen", "car")
    ,@(30, "car/car101.png", "passen", "car")
    ,@(31, "car/car099.png", "wehen", "car")
    ,@(32, "face/face109.png", "bergen", "face")
    ,@(33, "car/car096.png", "stoßen", "car")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
